# Standardize measurement values in the planet info spreadsheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Mass values: unify exponent to 10^25 ---
$ws.Range("F2").Value = ".033 × 10^25 kg"        # Mercury   was 3.285 × 10^23 kg
$ws.Range("F3").Value = ".487 × 10^25 kg"        # Venus     was 4.867 × 10^24 kg
$ws.Range("F4").Value = ".597 × 10^25 kg"        # Earth     was 5.972 × 10^24 kg
$ws.Range("F5").Value = ".064 × 10^25 kg"        # Mars      was 6.39 × 10^23 kg
$ws.Range("F6").Value = "189.8 × 10^25 kg"       # Jupiter   was 1.898 × 10^27 kg
$ws.Range("F7").Value = "56.83 × 10^25 kg"       # Saturn    was 5.683 × 10^26 kg
$ws.Range("F9").Value = "10.24 × 10^25 kg"       # Neptune   was 1.024 × 10^26 kg

# --- Distance from Sun: unify wording/units ---
$ws.Range("D3").Value = "67.03 million mi"       # Venus     was 67.028 million mi
$ws.Range("D8").Value = "1,784 million mi"       # Uranus    was 1.784 billion mi
$ws.Range("D9").Value = "2,793 million mi"       # Neptune   was 2.793 billion mi

# D9 previously carried a distinct leftover font (Tahoma, no charset); re-apply the
# standard Tahoma font so it matches the rest of the sheet's formatting.
$ws.Range("D9").Font.Name = "Tahoma"

# --- Column width tweaks (columns widened to fit the updated text) ---
$ws.Columns.Item(4).ColumnWidth = 16.333333333333332   # D: distance_from_sun -> ~17.19
$ws.Columns.Item(6).ColumnWidth = 16.166666666666668   # F: mass              -> ~17.06
$ws.Columns.Item(7).ColumnWidth = 13.166666666666666   # G: length_of_day     -> ~14.03

# --- Selection moved to C10 ---
[void]$ws.Range("C10").Select()
